$wb = $excel.ActiveWorkbook
$flags = $wb.Worksheets.Item(1)   # "Flags" sheet
$tests = $wb.Worksheets.Item(2)   # "Tests" sheet

# --- Flags sheet -----------------------------------------------------
# "Categories" flag value is cleared out.
$flags.Range("B3").ClearContents()

# "AllColors" flag flips from "False" to "True" (kept as literal text,
# not a boolean, since the column is Text-formatted). Compose the text
# in a scratch cell first (quote-prefixed so it isn't auto-coerced to a
# boolean), copy just the value across, then scrub the scratch cell so
# it leaves no trace in the sheet's used range.
$flags.Range("H1").Value = "'True"
$flags.Range("H1").Copy()
$flags.Range("B4").PasteSpecial(-4163)
$flags.Range("H1").Clear()

# Row 4 grows tall enough to show the wrapped "AllColors" description
# on two lines.
$flags.Rows.Item(4).RowHeight = 25.5

# Selection on the Flags sheet moves from B5 to A5.
$flags.Range("A5").Select()

# --- Tests sheet -------------------------------------------------------
# Add a new test case in row 42: format "h:m" applied to the same date
# serial used throughout the sheet, categorised under "Time".
$tests.Activate()
$tests.Range("B42").Value = "h:m"
$tests.Range("C42").Value = 17816.607951388887
$tests.Range("D42").Value = "Time"

# Selection on the Tests sheet moves from J18 to A43.
$tests.Range("A43").Select()
